$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Andre" row (placeholder/duplicate entry merged into Evelyn Cruz Mendoza)
$ws.Rows.Item(2).Delete()

# Fix name typos/corrections
$ws.Range("A11").Value = "Evelyn Cruz Mendoza"
$ws.Range("A29").Value = "Alejandro Rabelo García"

# Update "Modulos completados" (C) and "Modulo actual" (D) for each fellow
$ws.Range("C2").Value = 18
$ws.Range("D2").Value = "Conceptos básicos en torno a la incidencia en políticas públicas"
$ws.Range("C3").Value = 13
$ws.Range("D3").Value = "GOBIERNO ABIERTO Y LA AGENDA 2030 PARA EL DESARROLLO SOSTENIBLE"
$ws.Range("C4").Value = 15
$ws.Range("D4").Value = "Gobierno abierto como herramienta y plataforma para la incidencia (segunda parte); Metodología de Marco Lógico (evaluación final)"
$ws.Range("C5").Value = 18
$ws.Range("D5").Value = "Conceptos básicos en torno a la incidencia en políticas públicas"
$ws.Range("C6").Value = 15
$ws.Range("D6").Value = "Gobierno abierto como herramienta y plataforma para la incidencia (segunda parte); Metodología de Marco Lógico (evaluación final)"
$ws.Range("C7").Value = 14
$ws.Range("D7").Value = "Gobierno abierto como herramienta y plataforma para la incidencia (primera parte); Metodología de Marco Lógico (séptima parte)"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = "Gobierno abierto como herramienta y plataforma para la incidencia (primera parte); Metodología de Marco Lógico (séptima parte)"
$ws.Range("C9").Value = 13
$ws.Range("D9").Value = "GOBIERNO ABIERTO Y LA AGENDA 2030 PARA EL DESARROLLO SOSTENIBLE"
$ws.Range("C10").Value = 18
$ws.Range("D10").Value = "Conceptos básicos en torno a la incidencia en políticas públicas"
$ws.Range("C11").Value = 13
$ws.Range("D11").Value = "GOBIERNO ABIERTO Y LA AGENDA 2030 PARA EL DESARROLLO SOSTENIBLE"
$ws.Range("C12").Value = 15
$ws.Range("D12").Value = "Gobierno abierto como herramienta y plataforma para la incidencia (segunda parte); Metodología de Marco Lógico (evaluación final)"
$ws.Range("C13").Value = 12
$ws.Range("D13").Value = "Gobierno Abierto en la práctica mexicana (segunda parte)"
$ws.Range("C14").Value = 18
$ws.Range("D14").Value = "Conceptos básicos en torno a la incidencia en políticas públicas"
$ws.Range("C15").Value = 18
$ws.Range("D15").Value = "Conceptos básicos en torno a la incidencia en políticas públicas"
$ws.Range("C16").Value = 18
$ws.Range("D16").Value = "Conceptos básicos en torno a la incidencia en políticas públicas"
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = "Desarrollo Sostenible: Localizando los ODS (primera parte); Matriz de Marco Lógico (cuarta parte)"
$ws.Range("C18").Value = 18
$ws.Range("D18").Value = "Conceptos básicos en torno a la incidencia en políticas públicas"
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = "Gobierno Abierto como concepto (primera parte); metodología de Marco Lógico (quinta parte)"
$ws.Range("C20").Value = 18
$ws.Range("D20").Value = "Conceptos básicos en torno a la incidencia en políticas públicas"
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = "Gobierno abierto como herramienta y plataforma para la incidencia (primera parte); Metodología de Marco Lógico (séptima parte)"
$ws.Range("C22").Value = 15
$ws.Range("D22").Value = "Gobierno abierto como herramienta y plataforma para la incidencia (segunda parte); Metodología de Marco Lógico (evaluación final)"
$ws.Range("C23").Value = 18
$ws.Range("D23").Value = "Conceptos básicos en torno a la incidencia en políticas públicas"
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = "Gobierno abierto como herramienta y plataforma para la incidencia (segunda parte); Metodología de Marco Lógico (evaluación final)"
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = "Conceptos básicos en torno a la incidencia en políticas públicas"
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = "Conceptos básicos en torno a la incidencia en políticas públicas"
$ws.Range("C27").Value = 14
$ws.Range("D27").Value = "Gobierno abierto como herramienta y plataforma para la incidencia (primera parte); Metodología de Marco Lógico (séptima parte)"
$ws.Range("C28").Value = 14
$ws.Range("D28").Value = "Gobierno abierto como herramienta y plataforma para la incidencia (primera parte); Metodología de Marco Lógico (séptima parte)"
$ws.Range("C29").Value = 18
$ws.Range("D29").Value = "Conceptos básicos en torno a la incidencia en políticas públicas"
$ws.Range("C30").Value = 14
$ws.Range("D30").Value = "Gobierno abierto como herramienta y plataforma para la incidencia (primera parte); Metodología de Marco Lógico (séptima parte)"
$ws.Range("C31").Value = 13
$ws.Range("D31").Value = "GOBIERNO ABIERTO Y LA AGENDA 2030 PARA EL DESARROLLO SOSTENIBLE"
$ws.Range("C32").Value = 18
$ws.Range("D32").Value = "Conceptos básicos en torno a la incidencia en políticas públicas"
$ws.Range("C33").Value = 14
$ws.Range("D33").Value = "Gobierno abierto como herramienta y plataforma para la incidencia (primera parte); Metodología de Marco Lógico (séptima parte)"
$ws.Range("C34").Value = 15
$ws.Range("D34").Value = "Gobierno abierto como herramienta y plataforma para la incidencia (segunda parte); Metodología de Marco Lógico (evaluación final)"
$ws.Range("C35").Value = 5
$ws.Range("D35").Value = "Desarrollo Sostenible: ODS 16 (primera parte); Metodología de Marco Lógico (tercera parte)"
$ws.Range("C36").Value = 15
$ws.Range("D36").Value = "Gobierno abierto como herramienta y plataforma para la incidencia (segunda parte); Metodología de Marco Lógico (evaluación final)"
$ws.Range("C37").Value = 18
$ws.Range("D37").Value = "Conceptos básicos en torno a la incidencia en políticas públicas"
$ws.Range("C38").Value = 14
$ws.Range("D38").Value = "Gobierno abierto como herramienta y plataforma para la incidencia (primera parte); Metodología de Marco Lógico (séptima parte)"
$ws.Range("C39").Value = 15
$ws.Range("D39").Value = "Gobierno abierto como herramienta y plataforma para la incidencia (segunda parte); Metodología de Marco Lógico (evaluación final)"
$ws.Range("C40").Value = 18
$ws.Range("D40").Value = "Conceptos básicos en torno a la incidencia en políticas públicas"
$ws.Range("C41").Value = 15
$ws.Range("D41").Value = "Gobierno abierto como herramienta y plataforma para la incidencia (segunda parte); Metodología de Marco Lógico (evaluación final)"

# Update column D width to fit new (longer) content
$ws.Columns.Item(4).ColumnWidth = 153.248291